# "Drop in RMI script results for 3.0"
#
# 1. Remove the "Texas Data" sheet entirely (its commentary notes about the
#    NREL waste-heat correction go away with it).
# 2. Revert the HPEbP "natural gas reforming" efficiency formula (cell B3)
#    back to including waste heat in the input-energy denominator, i.e. from
#    118/(162+2) to 118/(162+2+46). The shared formulas across the rest of
#    row 3 (C3:AI3) all reference B3 (directly or transitively) so they will
#    recalculate automatically.

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

$texasSheet = $wb.Worksheets.Item("Texas Data")
$null = $texasSheet.Delete()

$hpebp = $wb.Worksheets.Item("HPEbP")
$hpebp.Range("B3").Formula = "=118/(162+2+46)"

$excel.DisplayAlerts = $true
